$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-11 Saturday" "2025-10-12 Sunday"

Replace-Text "55×97=" "86×54="
Replace-Text "31×26=" "68×16="
Replace-Text "68×61=" "11×14="
Replace-Text "52×52=" "76×16="
Replace-Text "48×31=" "67×74="

Replace-Text "74×85=" "22×40="
Replace-Text "72×17=" "80×33="
Replace-Text "87×26=" "12×21="
Replace-Text "52×43=" "55×99="
Replace-Text "24×43=" "51×35="

Replace-Text "48×49=" "44×79="
Replace-Text "90×92=" "70×37="
Replace-Text "91×82=" "54×89="
Replace-Text "82×50=" "99×89="
Replace-Text "19×11=" "24×62="

Replace-Text "70×49=" "47×95="
Replace-Text "87×70=" "87×24="
Replace-Text "14×78=" "26×12="
Replace-Text "24×58=" "41×35="
Replace-Text "38×99=" "97×76="

Replace-Text "21×99=" "26×17="
Replace-Text "47×50=" "63×45="
Replace-Text "25×51=" "79×18="
Replace-Text "97×58=" "20×83="
Replace-Text "29×83=" "79×11="
